# Update test data for naive and enhanced
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New trial measurements (Trial 1 .. Trial 10) for the "naive" rows (8-12)
# and the "enhanced" rows (14-18). Row 19 (n=1,000,000) trials are removed.
# NB: values are written in plain decimal (no E-notation) since the
# PowerShell parser here does not accept scientific-notation literals.
$data = @{
    8  = @(0.0000464916, 0.0000469685, 0.0000429153, 0.0000247955, 0.0000240803, 0.0000255108, 0.0000274181, 0.0000257492, 0.0000269413, 0.0000240803)
    9  = @(0.000191927, 0.0001952648, 0.0001897812, 0.0001943111, 0.0002081394, 0.0001881123, 0.0002028942, 0.0001897812, 0.0001835823, 0.0002529621)
    10 = @(0.0021913052, 0.0021038055, 0.002051115, 0.0020518303, 0.0020358562, 0.0020914078, 0.0020909309, 0.0020787716, 0.0021634102, 0.0021250248)
    11 = @(0.0251362324, 0.0265309811, 0.0254600048, 0.0249166489, 0.0249860287, 0.0255236626, 0.0251948833, 0.025187254, 0.0250411034, 0.0250582695)
    12 = @(0.3355109692, 0.325835228, 0.3291053772, 0.3261523247, 0.3278408051, 0.3202004433, 0.3270783424, 0.3341903687, 0.3430426121, 0.3250625134)
    14 = @(0.0000462532, 0.0000360012, 0.0000288486, 0.0000236034, 0.0000224113, 0.0000195503, 0.000020504, 0.0000259876, 0.0000219345, 0.0000207424)
    15 = @(0.000194788, 0.0001869202, 0.0001864433, 0.0001823902, 0.0001888275, 0.000180006, 0.0001823902, 0.000181675, 0.000171423, 0.0001797676)
    16 = @(0.0021557808, 0.0020956993, 0.0021114349, 0.0021185875, 0.0021226406, 0.0021467209, 0.0021729469, 0.0021443367, 0.0021195412, 0.0021443367)
    17 = @(0.0275919437, 0.0272498131, 0.0275504589, 0.0275506973, 0.0276648998, 0.0274918079, 0.0277087688, 0.027971983, 0.0279300213, 0.027859211)
    18 = @(0.4755814075, 0.4346909523, 0.466121912, 0.4473838806, 0.4455771446, 0.4395196438, 0.4684042931, 0.4553773403, 0.4605810642, 0.4459555149)
}

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# Row 19 (n=1,000,000) trial data was dropped entirely (cells removed, not
# just blanked); M19's AVERAGE formula goes with it, leaving the "points for
# graph" formula in N19 pointing at an empty M19.
$ws.Range("C19:L19").Clear()
$ws.Range("M19").ClearContents()

# Update the saved selection/active cell
$ws.Range("H19").Select()
